$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Total cost: sum of the individual room costs (D3:D7) into D8, formatted like them.
$ws.Range("D8").Formula = "=D3+D4+D5+D6+D7"
$ws.Range("D8").NumberFormat = $ws.Range("D3").NumberFormat

# Move the active selection to D9 (matches where the user clicked next).
$ws.Range("D9").Select()
